$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config")

# Shift the "index" values (column A) for rows 31..37 down by two positions,
# then fill in the newly vacated index values (27, 28) on rows 29 and 30.
# Working from the bottom up avoids any accidental overwrite collisions.
$ws.Range("A37").Value = 35
$ws.Range("A36").Value = 34
$ws.Range("A35").Value = 33
$ws.Range("A34").Value = 32
$ws.Range("A33").Value = 31
$ws.Range("A32").Value = 30
$ws.Range("A31").Value = 29
$ws.Range("A30").Value = 28
$ws.Range("A29").Value = 27

# Update the view: scroll so row 6 is the top row, and select K14.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K14").Select()
